$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension/measure labels in row 2
$ws.Range("G2").Value = "iaest-measure:sector-descripcion"
$ws.Range("H2").Value = "sdmx-dimension:refArea"
$ws.Range("I2").Value = "iaest-measure:sexo"

# Update row 3 (medida/dim) values
$ws.Range("G3").Value = "medida"
$ws.Range("I3").Value = "medida"

# Update row 4 values
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "URI-Comunidad"
$ws.Range("I4").Value = "xsd:int"

# Remove row 5 entirely (mapping-*.xlsx references)
$ws.Range("A5:K5").EntireRow.Delete()
